$p = $ppt.ActivePresentation

# Slide 6: remove "Group 4", "Group 6", and "TextBox 8" (Our Team heading),
# keeping only the decorative sidebar group ("Group 2").
$s6 = $p.Slides.Item(6)
for ($i = $s6.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.Name -ne "Group 2") {
        $shp.Delete()
    }
}

# Slide 7: remove all the "Name"/Lorem-ipsum team textboxes, the "Phases"
# heading, and the page-number textbox, keeping only the decorative
# sidebar group ("Group 9").
$s7 = $p.Slides.Item(7)
for ($i = $s7.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s7.Shapes.Item($i)
    if ($shp.Name -ne "Group 9") {
        $shp.Delete()
    }
}
